# Fruta / hortaliza, semanal
# Refresh the weekly price-sampling rows (Fecha, Volumen, Precio minimo/maximo/
# promedio ponderado, Precio $/Kg) for "Bruselas (repollito)" with the latest
# data pull. Row order / other columns (Mercado, Region, Categoria, etc.) are
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => Fecha(D), Volumen(J), Precio minimo(K), Precio maximo(L),
#        Precio promedio ponderado(M), Precio $/Kg(P)
$rows = @(
    @{ Row = 2;  D = 44391; J = 160; K = 20000; L = 20000; M = 20000; P = 1333 },
    @{ Row = 3;  D = 44399; J = 150; K = 22000; L = 22000; M = 22000; P = 1467 },
    @{ Row = 4;  D = 44392; J = 220; K = 23000; L = 23000; M = 23000; P = 1533 },
    @{ Row = 5;  D = 44476; J = 220; K = 20000; L = 22000; M = 20909; P = 1394 },
    @{ Row = 6;  D = 44365; J = 580; K = 20000; L = 22000; M = 21103; P = 1407 },
    @{ Row = 7;  D = 44483; J = 220; K = 18000; L = 20000; M = 18909; P = 1261 },
    @{ Row = 8;  D = 44446; J = 150; K = 22000; L = 24000; M = 22667; P = 1511 },
    @{ Row = 9;  D = 44396; J = 130; K = 22000; L = 22000; M = 22000; P = 1467 },
    @{ Row = 10; D = 44400; J = 130; K = 24000; L = 24000; M = 24000; P = 1600 },
    @{ Row = 11; D = 44406; J = 400; K = 20000; L = 22000; M = 20850; P = 1390 },
    @{ Row = 12; D = 44435; J = 140; K = 21000; L = 23000; M = 21714; P = 1448 },
    @{ Row = 13; D = 44449; J = 220; K = 22000; L = 24000; M = 23091; P = 1539 },
    @{ Row = 14; D = 44398; J = 130; K = 20000; L = 20000; M = 20000; P = 1333 },
    @{ Row = 15; D = 44699; J = 150; K = 18000; L = 20000; M = 18667; P = 1244 },
    @{ Row = 16; D = 44453; J = 280; K = 20000; L = 22000; M = 21286; P = 1419 }
)

foreach ($r in $rows) {
    $ws.Range("D$($r.Row)").Value = $r.D
    $ws.Range("J$($r.Row)").Value = $r.J
    $ws.Range("K$($r.Row)").Value = $r.K
    $ws.Range("L$($r.Row)").Value = $r.L
    $ws.Range("M$($r.Row)").Value = $r.M
    $ws.Range("P$($r.Row)").Value = $r.P
}
